$wb = $excel.ActiveWorkbook

# The worksheet that holds the tab/query table is the first sheet in the
# workbook (internally named "startup"), not the second "Sheet1".
$ws = $wb.Worksheets.Item(1)

# New "StatQuery" (Programs/Studies/Cases/Samples/Files) replaces the old,
# removed stat query used by all three rows (Cases/Samples/Files tabs).
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Dalmatian']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Update the saved view state: scrolled down one row, zoomed to 70%, and the
# B3 cell selected.
$excel.ActiveWindow.Zoom = 70
$ws.Range("B3").Select() | Out-Null
